# Swap the data in rows 3 and 4 of the "Staff" sheet.
# Before: row3 = ALANA  / DANSKIN / 3333 / MANAGER
#         row4 = NATHAN / DANSKIN / 2222 / ADMIN
# After:  row3 = NATHAN / DANSKIN / 2222 / ADMIN
#         row4 = ALANA  / DANSKIN / 3333 / MANAGER

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Staff")

# Read the original values. Value2 is used for reading because the plain
# Value property is exposed as an indexed/parameterized property that this
# PowerShell shim can't bind to directly (Value(Variant) {get}{set}).
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2

$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$c4 = $ws.Range("C4").Value2
$d4 = $ws.Range("D4").Value2

# The CODE column values ("3333" / "2222") are stored as text, not numbers.
# Writing a digit-only string straight back through Value2 makes Excel
# auto-convert it to a numeric cell. To keep it text (and keep the cell
# style untouched) write it as a quoted-string formula and then convert
# the formula to a static value via copy / paste-special values.
function Set-TextValue($range, [string]$value) {
    $range.Formula = '="' + $value + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Row 3 becomes the old row 4's data.
$ws.Range("A3").Value2 = $a4
$ws.Range("B3").Value2 = $b4
Set-TextValue $ws.Range("C3") $c4
$ws.Range("D3").Value2 = $d4

# Row 4 becomes the old row 3's data.
$ws.Range("A4").Value2 = $a3
$ws.Range("B4").Value2 = $b3
Set-TextValue $ws.Range("C4") $c3
$ws.Range("D4").Value2 = $d3

$excel.CutCopyMode = 0
